$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix hardcoded year in the "HOMEROOM 2LUHUR" header text
$ws.Range("B8").Value = "HOMEROOM 2LUHUR 2023"

# Update transaction values for row 16
$ws.Range("D16").Value = 3200
$ws.Range("E16").Value = 150

# Update transaction values for row 17
$ws.Range("D17").Value = 5670
$ws.Range("E17").Value = 710

# Update transaction values for row 36
$ws.Range("D36").Value = 100
